$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 value (was "23xddsss", now "23xddsss43")
$ws.Range("A2").Value = "23xddsss43"

# Remove row 3 entirely (was dsa3lx / QAZxsw124 / Automated / Test)
$ws.Rows.Item(3).Delete()

# Update the selected cell in the sheet view
$ws.Range("B5").Select()
